$d = $word.ActiveDocument

# Step 1: Justify all existing paragraphs (wdAlignParagraphJustify = 3)
foreach ($p in $d.Paragraphs) {
    $p.Alignment = 3
}

# Helper: find the 1-based index of the first paragraph whose text contains $needle.
function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# Step 2: Split "Descrizione delle operazioni di scambio client-server" paragraph
#   so that "client-server" is wrapped in proofErr spellcheck tags (text unchanged).
$idxDescr = Find-ParagraphIndex $d "Descrizione delle operazioni di scambio"
$paraDescr = $d.Paragraphs.Item($idxDescr)
$rngDescr = $paraDescr.Range
$xmlDescr = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Descrizione delle operazioni di scambio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>client-server</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
[void]$rngDescr.InsertXML($xmlDescr)

# Step 3: Rename heading "RICERCA RISTORANTE" -> "RECUPERO ITINERARIO"
#   (this is the "7-..." heading)
#   Use InsertXML so the existing run split (7 / - / RECUPERO ITINERARIO / (nuovo)) is preserved.
$idxHeading = Find-ParagraphIndex $d "7-RICERCA RISTORANTE"
$paraHeading = $d.Paragraphs.Item($idxHeading)
$rngHeading = $paraHeading.Range
$xmlHeading = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:ind w:left="708"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RECUPERO ITINERARIO</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> (nuovo)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
[void]$rngHeading.InsertXML($xmlHeading)

# Step 4: Replace the "Client invia al server..." paragraph with:
#   - the new descriptive paragraph about map retrieval
#   - an empty bold paragraph
#   - a new "8-RICERCA RISTORANTE (nuovo)" heading
#   - the original "Client invia al server..." paragraph text (moved down)
$idxClient = Find-ParagraphIndex $d "Client invia al server una richiesta contenente delle coordinate"
$paraClient = $d.Paragraphs.Item($idxClient)
$rngClient = $paraClient.Range
$xmlClient = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:ind w:left="708"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Quando l&#8217;utente seleziona una mappa tra quelle gi&#224; create in precedenza, il client invia una richiesta al server contenente il nome dell&#8217;utente e il nominativo della mappa; il server esegue una query sulla tabella delle mappe, restituendo il risultato ottenuto: una singola mappa, dato che nome utente e nominativo mappa sono le chiavi primarie della tabella.</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:left="708"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:ind w:left="708"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>8</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-RICERCA RISTORANTE</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> (nuovo)</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:left="708"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Client invia al server una richiesta contenente delle coordinate (la posizione prevista dall&#8217;itinerario nell&#8217;ora impostata per la sosta pasto) e un numero massimo di ristoranti da restituire;</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
[void]$rngClient.InsertXML($xmlClient)

Write-Host "Done"
